$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Columns.Item(5).ColumnWidth = 16.8
$ws.Range("E2").Value = "# FAMOUS CIGARS"
$ws.Range("E3").Value = "1502 CIGARS"
$ws.Range("F3").Value = "Years of blending, generations of expertise and centuries of tradition have all factored into the making of great cigars — and 1502 Cigars is no exception. Innovation, passion and development by companies such as theirs have been key in the effort to bring you the best cigar smoking experience known to man since he first rolled up a leaf and smoked it. And it’s hard to argue that with hard work like this, the premium cigar world is now the best it’s ever been. They’re committed to delighting the senses each and every time you set flame to the foot, with your preferences in mind. And that’s the beauty of all of the premium cigars on sale at Famous Smoke Shop — fitting all ranges of tastes, likes and prices, there is truly a cigar for everyone here at Famous.`nYou share our passion for a good cigar; and here, like every smoke we sell, when you buy 1502 Cigars online from Famous Smoke Shop you’re assured a flavorful experience, backed by our freshness guarantee. Buy 1502 Cigars on sale now, with confidence — and you’ll be telling your friends, “Here’s my new favorite thing from Famous Smoke Shop…”"
$ws.Rows.Item(3).EntireRow.AutoFit()
